# Краткий отчет — append 62 new product rows (GTIN + name) to the sheet,
# right after the current last data row (118), matching the existing
# table layout (column A = GTIN text, column B = full product name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new rows: push the far-below formatting-only rows
# (e.g. the stray styled cell at row 3449) down by 62 rows.
$ws.Rows("119:180").Insert()

$newRows = @(
    ,('4630460961793', 'Наволочка из тенселя 50 х 70 - М - Халкидики')
    ,('4630460961786', 'КПБ Постельное бельё полисатин Евро макси на резинке 180х200х30 - М - 80 - Витражи с наволочками 50х70')
    ,('4630460961779', 'Простыня из сатина-жаккарда 240 х 260 - М - Палладион')
    ,('4630460961762', 'Пододеяльник микросатин 240х260 - "Небосвод"')
    ,('4630460961755', 'Простыня из тенселя 200 х 220 - М - Королевский пурпур')
    ,('4630460961748', 'Пододеяльник из тенселя 200 х 220 - М - Пинклайт')
    ,('4630460961731', 'Наволочка из тенселя 70 х 70 - М - Хрустальный оазис')
    ,('4630460961724', 'Пододеяльник из тенселя 175 х 215 - М - Гималайские маки')
    ,('4630460961717', 'Простыня микросатин 220х240 - "Метель"')
    ,('4630460961700', 'Простыня круглая из тенселя на резинке 220х220х40 - М - 60S Белый')
    ,('4630460961694', 'Пододеяльник из сатина 145х215 - "Небесный глянцевая"')
    ,('4630460961687', 'Пододеяльник из сатина 145 х 215 - М - 125 - Голубая бирюза')
    ,('4630460961670', 'Наволочка из сатина 50 х 70 - М - 130 - Ночная синева')
    ,('4630460961663', 'Простыня Жатка 180х220 - "Парижский Шарм"')
    ,('4630460961656', 'Пододеяльник твил-сатин 150 х 200 - "Дино-пати"')
    ,('4630460961649', 'КПБ Постельное бельё тенсель семейное на резинке 180х200х30 - М - 60S 129 Темно-зеленый с наволочками 50х70')
    ,('4630460961632', 'Простыня из бязи на резинке 80 х 200 х 30 - М - 142 ГОСТ - Спящий Мишка на бежевом')
    ,('4630460961625', 'Простыня из страйп-сатина 260 х 260 - М - 135 - Эмерэльд 1х1')
    ,('4630460961618', 'Пододеяльник тенсель-жаккард 200 х 220 - "Табия"')
    ,('4630460961601', 'Наволочка тенсель 50х70 - "Коньячный"')
    ,('4630460961595', 'Наволочка тенсель 50х50 - "Коньячный"')
    ,('4630460961588', 'Простыня Микрофибра 300х300 - "Эстрелья"')
    ,('4630460961571', 'Простыня из страйп-сатина 280 х 280 - М - 135 - Белый 1х1')
    ,('4630460961564', 'Наволочка Жатка 50х70 - "Снежный Хищник"')
    ,('4630460961557', 'Наволочка из бязи 70 х 70 - М - 125 - Сансет')
    ,('4630460961540', 'Простыня из сатина на резинке 200 х 200 х 40 - М - 125 - Голубая бирюза')
    ,('4630460961533', 'Наволочка из тенселя 40 х 60 - М - 60S 4292 Кремовый')
    ,('4630460961526', 'КПБ Постельное бельё из страйп-сатина семейное на резинке 160х200х30 - М - 130 - Пел Догвуд 1х1')
    ,('4630460961519', 'Простыня Жатка 300х300 - "Леопардовый Сон"')
    ,('4630460961502', 'Пододеяльник Жатка 220х240 - "Леопардовый Сон"')
    ,('4630460961496', 'Простыня из страйп-сатина на резинке 200 х 220 х 40 - М - 135 - Углерод 1х1')
    ,('4630460961489', 'Простыня из страйп-сатина 240 х 260 - М - 125 - Шампань 3х3')
    ,('4630460961472', 'Наволочка из сатина 50 х 70 - М - Изольда')
    ,('4630460961465', 'Простыня микросатин 180х220 - "Легкий роман"')
    ,('4630460961458', 'Простыня круглая из сатина на резинке 230 х 230 х 30 - М - 125 - Вишневые лепестки')
    ,('4630460961441', 'Пододеяльник страйп-сатин 175х215 - 125 - "Миндаль" 1х1')
    ,('4630460961434', 'Простыня круглая из перкаля на резинке 200 х 200 х 30 - М - 110 - Каори на зеленом')
    ,('4630460961427', 'Пододеяльник из поплина 200 х 200 - М - Феста А+В')
    ,('4630460961410', 'Наволочка из перкаля 70 х 70 - Р - М - Дуэт')
    ,('4630460961403', 'Наволочка поплин 70 х 70 - "Оллорины"')
    ,('4630460961397', 'КПБ Постельное бельё тенсель-жаккард 2 спальное на резинке 160х200х30 - "Табия" с наволочками 50х70')
    ,('4630460961380', 'Пододеяльник из страйп микрофибры 220 х 240 - М - 85 - Белый 3х3')
    ,('4630460961373', 'Пододеяльник из сатина 200 х 220 - М - Красный')
    ,('4630460961366', 'Наволочка тенсель 40 х 40 - "Белый Страйп"')
    ,('4630460961359', 'КПБ Постельное белье тенсель 1,5 спальное - М - Черный блеск с наволочками 50х70')
    ,('4630460961342', 'КПБ Постельное бельё страйп-сатин Евро макси на резинке 180х200х30 - М - 135 - Марриотт Белый 3х3')
    ,('4630460961335', 'Простыня из сатина 260 х 260 - М - 140 - Можиана Люкс')
    ,('4630460961328', 'Простыня из перкаля на резинке 120 х 200 х 30 - М - 125 - Лунный песок')
    ,('4630460961311', 'Простыня тенсель на резинке 200х220х20 - "Эрос"')
    ,('4630460961304', 'Наволочка страйп-сатин 50х70 - 125 - "Индиго" 1х1')
    ,('4630460961298', 'Пододеяльник варёный хлопок 200х220 - "Стальная тень" (меланж)')
    ,('4630460961281', 'Простыня страйп-сатин на резинке 160х200х30 - 125 - "Олива" 1х1')
    ,('4630460961274', 'КПБ Постельное бельё тенсель Евро на резинке 160х200х30 - М - Блу Скай')
    ,('4630460961267', 'КПБ Постельное бельё полисатин семейное - "Глубь"')
    ,('4630460961250', 'КПБ Постельное бельё Жатка Евро на резинке 160х200х30 - "Замороженный Перец"')
    ,('4630460961243', 'Простыня из поплина 220 х 240 - М - Гранола')
    ,('4630460961236', 'Простыня полисатин 300х300 - "Синиль"')
    ,('4630460961229', 'КПБ Постельное бельё ранфорс 1,5 спальное на резинке 140х200х30 - "Фантасмагория" с наволочками 50х70')
    ,('4630460961212', 'Простыня Жатка 180х220 - "Ромашковая Нежность"')
    ,('4630460961205', 'КПБ Постельное бельё перкаль Евро - М - 110 - Каори на зеленом с наволочками 50х70')
    ,('4630460961199', 'Пододеяльник из тенселя 150 х 200 - М - Пинклайт')
    ,('4630460961182', 'Пододеяльник тенсель 150 х 200 - "Серый Райнар"')
)

$r = 119
foreach ($pair in $newRows) {
    $gtinCell = $ws.Cells.Item($r, 1)
    # Leading apostrophe forces the numeric-looking GTIN to be stored as
    # text (matching the shared-string type used throughout column A);
    # resetting the style afterwards avoids leaving a quote-prefix format
    # on the cell so it stays visually identical to the rest of the column.
    $gtinCell.Value = "'" + $pair[0]
    $gtinCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r++
}
